$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 17; this shifts existing rows 17:86 down to 18:87
$ws.Rows("17:17").Insert()

# Fill the new row 17 with the new weekly record.
# Columns A, B, C, E, F, G, H, I, R are constant across all data rows in this sheet.
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44672
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 100112030
$ws.Range("G17").Value = "Poroto granado"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 360
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 17500
$ws.Range("N17").Value = "$/malla 25 kilos"
$ws.Range("O17").Value = "Provincia de Limarí"
$ws.Range("P17").Value = 700
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
